$d = $word.ActiveDocument

# 1) "Entrega: 22-05-2013" -> "Entrega: 15-05-2013"
#    Only the date-number run (" 22") changes; the preceding bold
#    "Entrega:" run is left alone.
$r = $d.Content
$found = $r.Find.Execute("Entrega: 22-05-2013", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $r.Start
    $sub = $d.Range($start + 8, $start + 11)
    $sub.Text = " 15"
}

# 2) Footer page-number field cached result "8" -> "6"
$footer = $d.Sections(1).Footers(1)
$ch = $footer.Range.Characters(1)
if ($ch.Text -eq "8") {
    $ch.Text = "6"
}
